$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update ASV_rank (column G) values from 42 to 41 for rows 9-12
$ws.Range("G9").Value = 41
$ws.Range("G10").Value = 41
$ws.Range("G11").Value = 41
$ws.Range("G12").Value = 41
